$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Pthlh"
$ws.Range("C2").Value = "Pth1r"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2756003333333333
$ws.Range("H2").Value = 0.8268009999999999
$ws.Range("I2").Value = 0.03520863368480177
$ws.Range("J2").Value = 0.03520863368480177
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.9990806666666666
$ws.Range("N2").Value = 2.997242
$ws.Range("O2").Value = 0.1356546064507813
$ws.Range("P2").Value = 0.1356546064507813
$ws.Range("Q2").Value = 0.2753469647602221
$ws.Range("R2").Value = 2.478122682842
$ws.Range("S2").Value = 0.004776213346181504
$ws.Range("T2").Value = 0.004776213346181507

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pthlh"
$ws.Range("C3").Value = "Pth1r"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2756003333333333
$ws.Range("H3").Value = 0.8268009999999999
$ws.Range("I3").Value = 0.03520863368480177
$ws.Range("J3").Value = 0.03520863368480177
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.268092
$ws.Range("N3").Value = 12.804276
$ws.Range("O3").Value = 0.5795191117925025
$ws.Range("P3").Value = 0.5795191117925026
$ws.Range("Q3").Value = 1.176287577897333
$ws.Range("R3").Value = 10.586588201076
$ws.Range("S3").Value = 0.02040407612044391
$ws.Range("T3").Value = 0.02040407612044391

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Pthlh"
$ws.Range("C4").Value = "Pth1r"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2756003333333333
$ws.Range("H4").Value = 0.8268009999999999
$ws.Range("I4").Value = 0.03520863368480177
$ws.Range("J4").Value = 0.03520863368480177
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.318589
$ws.Range("N4").Value = 0.955767
$ws.Range("O4").Value = 0.04325783378307253
$ws.Range("P4").Value = 0.04325783378307254
$ws.Range("Q4").Value = 0.08780323459633332
$ws.Range("R4").Value = 0.7902291113669999
$ws.Range("S4").Value = 0.001523049223666243
$ws.Range("T4").Value = 0.001523049223666244

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Pthlh"
$ws.Range("C5").Value = "Pth1r"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2756003333333333
$ws.Range("H5").Value = 0.8268009999999999
$ws.Range("I5").Value = 0.03520863368480177
$ws.Range("J5").Value = 0.03520863368480177
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.779124
$ws.Range("N5").Value = 5.337372
$ws.Range("O5").Value = 0.2415684479736436
$ws.Range("P5").Value = 0.2415684479736436
$ws.Range("Q5").Value = 0.4903271674413333
$ws.Range("R5").Value = 4.412944506972
$ws.Range("S5").Value = 0.00850529499451011
$ws.Range("T5").Value = 0.008505294994510113

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Pthlh"
$ws.Range("C6").Value = "Pth1r"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 7.552034666666667
$ws.Range("H6").Value = 22.656104
$ws.Range("I6").Value = 0.9647913663151982
$ws.Range("J6").Value = 0.9647913663151982
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.9990806666666666
$ws.Range("N6").Value = 2.997242
$ws.Range("O6").Value = 0.1356546064507813
$ws.Range("P6").Value = 0.1356546064507813
$ws.Range("Q6").Value = 7.545091829463111
$ws.Range("R6").Value = 67.90582646516799
$ws.Range("S6").Value = 0.1308783931045998
$ws.Range("T6").Value = 0.1308783931045998

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Pthlh"
$ws.Range("C7").Value = "Pth1r"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 7.552034666666667
$ws.Range("H7").Value = 22.656104
$ws.Range("I7").Value = 0.9647913663151982
$ws.Range("J7").Value = 0.9647913663151982
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.268092
$ws.Range("N7").Value = 12.804276
$ws.Range("O7").Value = 0.5795191117925025
$ws.Range("P7").Value = 0.5795191117925026
$ws.Range("Q7").Value = 32.23277874452267
$ws.Range("R7").Value = 290.095008700704
$ws.Range("S7").Value = 0.5591150356720586
$ws.Range("T7").Value = 0.5591150356720587

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Pthlh"
$ws.Range("C8").Value = "Pth1r"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 7.552034666666667
$ws.Range("H8").Value = 22.656104
$ws.Range("I8").Value = 0.9647913663151982
$ws.Range("J8").Value = 0.9647913663151982
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.318589
$ws.Range("N8").Value = 0.955767
$ws.Range("O8").Value = 0.04325783378307253
$ws.Range("P8").Value = 0.04325783378307254
$ws.Range("Q8").Value = 2.405995172418667
$ws.Range("R8").Value = 21.653956551768
$ws.Range("S8").Value = 0.04173478455940628
$ws.Range("T8").Value = 0.0417347845594063

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Pthlh"
$ws.Range("C9").Value = "Pth1r"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 7.552034666666667
$ws.Range("H9").Value = 22.656104
$ws.Range("I9").Value = 0.9647913663151982
$ws.Range("J9").Value = 0.9647913663151982
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.779124
$ws.Range("N9").Value = 5.337372
$ws.Range("O9").Value = 0.2415684479736436
$ws.Range("P9").Value = 0.2415684479736436
$ws.Range("Q9").Value = 13.43600612429867
$ws.Range("R9").Value = 120.924055118688
$ws.Range("S9").Value = 0.2330631529791335
$ws.Range("T9").Value = 0.2330631529791335
